$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.793.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "'3.228.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'577.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'173.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.49%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'3.227.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "'0.390"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("D13").Value = "'3.794.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "'64.939.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "'25.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.247.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("B18").Value = "'ShibaInu"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'0.0000159"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").Value = "'413.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.47%  "
$ws.Range("D20").Value = "'5.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").Value = "'12.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "'7.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'70.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("D27").Value = "'0.493"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "'0.0000110"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("D29").Value = "'9.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'1.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("D32").Value = "'21.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "'4.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("D35").Value = "'6.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").Value = "'156.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").Value = "'1.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "'2.832.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("D40").Value = "'1.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").Value = "'25.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.43%  "
$ws.Range("D42").Value = "'4.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").Value = "'0.730"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.55%  "
$ws.Range("D44").Value = "'39.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").Value = "'5.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.54%  "
$ws.Range("D46").Value = "'0.0628"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("D47").Value = "'305.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.65%  "
$ws.Range("D48").Value = "'2.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.75%  "
$ws.Range("D49").Value = "'22.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.62%  "
$ws.Range("D50").Value = "'0.0263"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("E51").Value = "  -0.53%  "
